# Simulator full-month coverage, persist logs, fix employees
# Updates the "Weekly Timesheet" and "Jason Schema" sheets for Phil Henderson's
# 2026-01-19 week: corrects client/employee names, hours, rates and totals,
# and fixes the employee ID.

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: "Weekly Timesheet" ----
$ws1 = $wb.Worksheets.Item("Weekly Timesheet")

# Row 2: 2026-01-19 - PTO day
$ws1.Range("B2").Value = "PTO"
$ws1.Range("C2").Value = 6
$ws1.Range("D2").Value = "PTO"
$ws1.Range("E2").Value = 90
$ws1.Range("F2").Value = 540

# Row 3: 2026-01-20 - Keevil
$ws1.Range("B3").Value = "Keevil"
$ws1.Range("C3").Value = 6.5
$ws1.Range("E3").Value = 90
$ws1.Range("F3").Value = 585

# Row 4: 2026-01-21 - Howard
$ws1.Range("B4").Value = "Howard"
$ws1.Range("C4").Value = 6
$ws1.Range("E4").Value = 90
$ws1.Range("F4").Value = 540

# Row 5: 2026-01-22 - Markfield
$ws1.Range("B5").Value = "Markfield"
$ws1.Range("C5").Value = 6
$ws1.Range("E5").Value = 90
$ws1.Range("F5").Value = 540

# Row 6: 2026-01-23 - Layne
$ws1.Range("B6").Value = "Layne"
$ws1.Range("C6").Value = 7.5
$ws1.Range("E6").Value = 90
$ws1.Range("F6").Value = 675

# Row 8: SUBTOTAL
$ws1.Range("C8").Value = 32
$ws1.Range("D8").Value = "Reg: 32 / OT: 0"
$ws1.Range("F8").Value = 2880

# Row 11: HOURLY SUBTOTAL
$ws1.Range("F11").Value = 2880

# Row 13: GRAND TOTAL
$ws1.Range("F13").Value = 2880

# ---- Sheet 2: "Jason Schema" ----
$ws2 = $wb.Worksheets.Item("Jason Schema")

# Row 2: 2026-01-19 - PTO day
$ws2.Range("B2").Value = "emp_75yd72zj"
$ws2.Range("D2").Value = "PTO"
$ws2.Range("E2").Value = 6
$ws2.Range("F2").Value = 90
$ws2.Range("G2").Value = 540
$ws2.Range("H2").Value = "PTO"
$ws2.Range("I2").Value = "PTO"

# Row 3: 2026-01-20 - Keevil
$ws2.Range("B3").Value = "emp_75yd72zj"
$ws2.Range("D3").Value = "Keevil"
$ws2.Range("E3").Value = 6.5
$ws2.Range("F3").Value = 90
$ws2.Range("G3").Value = 585

# Row 4: 2026-01-21 - Howard
$ws2.Range("B4").Value = "emp_75yd72zj"
$ws2.Range("D4").Value = "Howard"
$ws2.Range("E4").Value = 6
$ws2.Range("F4").Value = 90
$ws2.Range("G4").Value = 540

# Row 5: 2026-01-22 - Markfield
$ws2.Range("B5").Value = "emp_75yd72zj"
$ws2.Range("D5").Value = "Markfield"
$ws2.Range("E5").Value = 6
$ws2.Range("F5").Value = 90
$ws2.Range("G5").Value = 540

# Row 6: 2026-01-23 - Layne
$ws2.Range("B6").Value = "emp_75yd72zj"
$ws2.Range("D6").Value = "Layne"
$ws2.Range("E6").Value = 7.5
$ws2.Range("F6").Value = 90
$ws2.Range("G6").Value = 675

$wb.Save()
